$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CDW")

# Row 4 - Inventory
$ws.Range("B4").Value = 760000000.0
$ws.Range("C4").Value = 648000000.0
$ws.Range("D4").Value = 695000000.0
$ws.Range("E4").Value = 672000000.0
$ws.Range("F4").Value = 611000000.0

# Row 13 - Accounts Payable
$ws.Range("B13").Value = 2088000000.0
$ws.Range("C13").Value = 1922000000.0
$ws.Range("D13").Value = 1809000000.0
$ws.Range("E13").Value = 1967000000.0
$ws.Range("F13").Value = 1835000000.0

# Row 23 - Long Term Tax Liability (Deferred)
$ws.Range("B23").Value = 55000000.0
$ws.Range("C23").Value = 63000000.0
$ws.Range("D23").Value = 68000000.0
$ws.Range("E23").Value = 63000000.0
$ws.Range("F23").Value = 62000000.0

$wb.Save()
